# Auto-generated Excel COM-interop script
# Applies numeric updates to F (想去人数) and G (最低票价) columns
# across sheets 展览, 演出, 本地生活, 全部类型 per the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (37 cell updates) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 527
$ws.Range("F4").Value = 579
$ws.Range("F5").Value = 9229
$ws.Range("F7").Value = 11955
$ws.Range("G7").Value = 238
$ws.Range("F8").Value = 11955
$ws.Range("G8").Value = 238
$ws.Range("F12").Value = 122
$ws.Range("F16").Value = 2048
$ws.Range("F17").Value = 812
$ws.Range("F18").Value = 771
$ws.Range("F19").Value = 393
$ws.Range("F20").Value = 37
$ws.Range("F21").Value = 408
$ws.Range("F23").Value = 96
$ws.Range("F24").Value = 662
$ws.Range("F25").Value = 58
$ws.Range("F26").Value = 15
$ws.Range("F27").Value = 1525
$ws.Range("F29").Value = 20
$ws.Range("F30").Value = 17
$ws.Range("F33").Value = 1441
$ws.Range("F34").Value = 6
$ws.Range("F36").Value = 339
$ws.Range("F37").Value = 512
$ws.Range("F38").Value = 392
$ws.Range("F39").Value = 2151
$ws.Range("F41").Value = 73
$ws.Range("G41").Value = 39.9
$ws.Range("F42").Value = 155
$ws.Range("F43").Value = 574
$ws.Range("F44").Value = 443
$ws.Range("F45").Value = 155
$ws.Range("F46").Value = 880
$ws.Range("F47").Value = 678
$ws.Range("F49").Value = 286
$ws.Range("F50").Value = 257

# --- Sheet: 演出 (4 cell updates) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 71
$ws.Range("F18").Value = 84
$ws.Range("F23").Value = 77
$ws.Range("F25").Value = 65

# --- Sheet: 本地生活 (3 cell updates) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2906
$ws.Range("F4").Value = 368
$ws.Range("F6").Value = 253

# --- Sheet: 全部类型 (35 cell updates) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 527
$ws.Range("F5").Value = 368
$ws.Range("F7").Value = 579
$ws.Range("F8").Value = 9229
$ws.Range("F10").Value = 11955
$ws.Range("G10").Value = 238
$ws.Range("F12").Value = 122
$ws.Range("F14").Value = 2048
$ws.Range("F15").Value = 812
$ws.Range("F16").Value = 771
$ws.Range("F17").Value = 393
$ws.Range("F18").Value = 37
$ws.Range("F19").Value = 408
$ws.Range("F22").Value = 662
$ws.Range("F25").Value = 15
$ws.Range("F26").Value = 253
$ws.Range("F27").Value = 1525
$ws.Range("F31").Value = 71
$ws.Range("F32").Value = 1441
$ws.Range("F34").Value = 6
$ws.Range("F36").Value = 84
$ws.Range("F37").Value = 512
$ws.Range("F38").Value = 392
$ws.Range("F39").Value = 2151
$ws.Range("F40").Value = 73
$ws.Range("G40").Value = 39.9
$ws.Range("F41").Value = 155
$ws.Range("F42").Value = 574
$ws.Range("F43").Value = 443
$ws.Range("F44").Value = 155
$ws.Range("F45").Value = 77
$ws.Range("F46").Value = 65
$ws.Range("F48").Value = 678
$ws.Range("F49").Value = 286
$ws.Range("F50").Value = 258
